$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effort R 1.0")

# --- Fix existing rows ---

# Row 12: B12 1.5 -> 2, C12 2.5 -> 2
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 2

# Row 18: the 2.25 entry was recorded in the wrong column (C18); move it to B18
$ws.Range("B18").Value = 2.25
$ws.Range("C18").ClearContents()

# --- New entries: TC14 work started ---

# Row 21
$ws.Range("A21").Value = 41452
$ws.Range("B21").Value = 2.5
$ws.Range("D21").Value = "TODOs, code cleanup"

# Row 22
$ws.Range("A22").Value = 41455
$ws.Range("B22").Value = 1.5
$ws.Range("D22").Value = "Concept of tc14"

# Row 23
$ws.Range("A23").Value = 41456
$ws.Range("B23").Value = 2
$ws.Range("D23").Value = "Implementation tc14"

# Update the selection to the next empty row, matching the author's cursor position
$ws.Range("A24").Select()
